# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strike#) values for rows 2-20 in column G
$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 3
    12 = 1
    13 = 1
    14 = 1
    16 = 1
    17 = 2
    18 = 3
    19 = 1
    20 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
